$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45172 (2023-09-03) to 45175 (2023-09-06) for every data row
# (rows 2 through 398).
$ws.Range("C2:C398").Value = 45175
